# Apply "Add data for 2022-08-06" update to the carjacking-by-neighborhood-by-month
# workbook: bump the "through July 28" snapshot to "through July 29" (sheet name,
# header label in B1) and update the affected neighborhood/month cell counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet and update the column header text to reflect the new "as of" date.
$ws.Name = "Through 2022-07-29"
$ws.Range("B1").Value = "July 2022 (through July 29)"

# Updated cell counts (row = neighborhood, column = month "through day 29" snapshot).
$ws.Range("B2").Value = 15
$ws.Range("P2").Value = 8
$ws.Range("AR2").Value = 5

$ws.Range("I4").Value = 2

$ws.Range("B5").Value = 9

$ws.Range("B6").Value = 10
$ws.Range("AD6").Value = 3

$ws.Range("AK7").Value = 6

$ws.Range("P8").Value = 22

$ws.Range("AD12").Value = 2

$ws.Range("P18").Value = 2

$ws.Range("I23").Value = 1

$ws.Range("P26").Value = 3

$ws.Range("AR29").Value = 1

$ws.Range("W34").Value = 1

$ws.Range("AR35").Value = 1

$ws.Range("I47").Value = 2

$ws.Range("P52").Value = 10

$ws.Range("B53").Value = 1

$ws.Range("W54").Value = 1

$ws.Range("P57").Value = 2
$ws.Range("P58").Value = 1

$ws.Range("P61").Value = 1
